$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Data")
$ws2 = $wb.Worksheets.Item("Codebook")

# Recharacterize "Political Party" / "Political party" header text into "Party"
# on both the Data sheet and the Codebook sheet (reusing the same string).
$ws1.Range("E1").Value = "Party"
$ws2.Range("A6").Value = "Party"

# Update view/selection state:
# Codebook becomes unselected with its cursor left at B6,
# Data becomes the active/selected sheet with its cursor at F1.
$ws2.Range("B6").Select()
$ws1.Activate()
$ws1.Range("F1").Select()
